$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same style as the other header cells (e.g. H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J (rows 2-11)
$values = @{
    2  = @(8, 8)
    3  = @(13, 13)
    4  = @(9, 9)
    5  = @(5, 5)
    6  = @(7, 7)
    7  = @(6, 6)
    8  = @(6, 6)
    9  = @(5, 5)
    10 = @(4, 4)
    11 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
}
